# I2C master/slave comm; library skeletons
# BOM update: the "12V Motor" part is swapped out for a "6V Motor" sourced
# from Pololu (was RobotShop), with a new unit price / quantity, and the
# cursor selection moves to D9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the row-4 (Part) style before we touch anything, so we can put
# it back after the hyperlink swap (adding a hyperlink re-styles the cell).
$origStyle = $ws.Range("A4").Style

# Drop the old hyperlink that lived on A4 (it pointed at the RobotShop
# 12V-motor listing) - it no longer applies to the new part.
foreach ($hl in @($ws.Hyperlinks)) {
    if ($hl.Range.Row -eq 4 -and $hl.Range.Column -eq 1) {
        $hl.Delete()
    }
}

# New part data for row 4: 6V Motor from Pololu, $40/unit, qty 2.
$ws.Range("A4").Value = "6V Motor"
$ws.Range("B4").Value = "Pololu"
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = 2

# Point A4 at the Pololu/McMaster-style product page.
$ws.Hyperlinks.Add($ws.Range("A4"), "https://www.mcmaster.com/")

# Restore the original (non-hyperlink-wizard) cell style for A4.
$ws.Range("A4").Style = $origStyle

# Move the active selection to D9.
$ws.Range("D9").Select()
